$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices corresponding to B, C, D, E, F, I, J, K, L, M, N (column H is intentionally skipped, as in the source data)
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

# New vm_pu values for rows 2..25, one row per array, values ordered to match $cols
$data = @(
    @(1.02, 1.030763952664216, 1.038882333672708, 1.048491652779333, 1.052714059346016, 1.032864540483332, 1.035903398080549, 1.041669102157293, 1.051251348547466, 1.0554620285946, 1.037374499725604),
    @(1.02, 1.031697441858328, 1.039587469067481, 1.049423934783364, 1.053629788691581, 1.032992033668823, 1.036478560065953, 1.04218461924586, 1.051995375568136, 1.056190394774439, 1.037950478507),
    @(1.02, 1.032301953380738, 1.040043999985463, 1.050028405245716, 1.05422325528672, 1.033073247534495, 1.036850597200812, 1.042517780656738, 1.05247740363835, 1.056661992261293, 1.038323043977263),
    @(1.02, 1.032556204255269, 1.040235986605266, 1.050282815602648, 1.054472969390784, 1.033107082185269, 1.037006969194764, 1.042657741814616, 1.052680188892259, 1.05686032136182, 1.038479638037382),
    @(1.02, 1.032598900763562, 1.040268225538178, 1.050325549240513, 1.0545149103879, 1.033112745114364, 1.037033222809866, 1.042681236018247, 1.052714245639684, 1.056893625728846, 1.038505928935629),
    @(1.02, 1.032305350245844, 1.040046565082188, 1.050031803547899, 1.054226591113334, 1.033073700844347, 1.036852686779712, 1.04251965121852, 1.052480112713762, 1.056664642072934, 1.038325136523605),
    @(1.02, 1.031079329350076, 1.039120582501931, 1.048806467994306, 1.053023341197625, 1.032907892664203, 1.036097803301305, 1.041843408333976, 1.051502672446171, 1.055708121304672, 1.037569181024072),
    @(1.02, 1.028922659325426, 1.037490949402998, 1.04665669278526, 1.050910239808253, 1.032605918108548, 1.034766643112852, 1.040648671170614, 1.049784894501414, 1.054024933501259, 1.036236130435466),
    @(1.02, 1.027487453040998, 1.036406005697169, 1.04522993313064, 1.049506422781145, 1.03239804668426, 1.033878618595449, 1.039850154287971, 1.048642875162815, 1.052904448554767, 1.03534684482139),
    @(1.02, 1.026866616075928, 1.035936582366411, 1.044613672618995, 1.048899739755492, 1.03230648827043, 1.033493968113815, 1.039503921037398, 1.048149134855414, 1.052419671472246, 1.034961648091952),
    @(1.02, 1.026636103234618, 1.035762273880554, 1.044384998248645, 1.048674569309512, 1.032272247258678, 1.033351073249502, 1.039375245039126, 1.047965853396208, 1.052239665104109, 1.034818550300541),
    @(1.02, 1.026685544767932, 1.035799661073651, 1.044434039144975, 1.048722861026345, 1.032279602568138, 1.033381725506546, 1.039402849607552, 1.048005162636716, 1.052278274280856, 1.034849246087306),
    @(1.02, 1.026847559896439, 1.035922172816725, 1.044594765578362, 1.048881123443202, 1.032303662630198, 1.033482156752456, 1.039493286051824, 1.048133982385841, 1.052404790829211, 1.034949819957105),
    @(1.02, 1.02694739520323, 1.035997663841815, 1.044693825254705, 1.048978657858558, 1.03231845608376, 1.03354403332923, 1.039548997745933, 1.048213367850389, 1.052482750028111, 1.035011784405715),
    @(1.02, 1.027528669024982, 1.036437167583048, 1.04527086479582, 1.049546711304693, 1.032404090530786, 1.033904143925913, 1.039873122835155, 1.048675659233892, 1.052936630156573, 1.035372406400751),
    @(1.02, 1.027893452568064, 1.036712955291849, 1.045633238686051, 1.049903352943372, 1.032457392500487, 1.034129997676656, 1.040076312828025, 1.048965847191538, 1.053221445246759, 1.035598580889726),
    @(1.02, 1.028106284086042, 1.036873852661365, 1.045844753331627, 1.050111489633468, 1.032488333172983, 1.034261721738688, 1.0401947846757, 1.049135182317132, 1.053387611606394, 1.035730492015034),
    @(1.02, 1.028178864187472, 1.036928720408825, 1.045916899432715, 1.050182478098218, 1.032498857774617, 1.034306634071084, 1.040235172774189, 1.049192933630692, 1.053444276538233, 1.035775468128095),
    @(1.02, 1.027854308589979, 1.036683362243053, 1.045594344069136, 1.049865076903861, 1.032451689159529, 1.034105767003816, 1.040054517156721, 1.048934705168882, 1.053190883287765, 1.035574315806551),
    @(1.02, 1.026799847886121, 1.035886094597486, 1.044547429227344, 1.048834514162494, 1.032296583945493, 1.033452582757858, 1.039466656684342, 1.048096044977214, 1.052367533153101, 1.034920203964043),
    @(1.02, 1.026137408910127, 1.035385146908372, 1.043890536617709, 1.048187592781626, 1.032197720316111, 1.033041793067042, 1.039096643598522, 1.047569415908119, 1.051850216268309, 1.034508830904747),
    @(1.02, 1.026488528611199, 1.035650677348307, 1.044238639899443, 1.048530439500072, 1.03225025690141, 1.033259570252313, 1.039292832266922, 1.04784852802624, 1.052124421560383, 1.0347269173586),
    @(1.02, 1.027871995885117, 1.036696733959802, 1.045611918414916, 1.049882371847521, 1.032454266715684, 1.034116715840926, 1.040064365822899, 1.048948776681764, 1.053204692802568, 1.035585280192265),
    @(1.02, 1.02947976115684, 1.03791199455779, 1.047211335458557, 1.051455667116902, 1.032685143597456, 1.035110887359405, 1.040957901202036, 1.050228428779765, 1.054459794483524, 1.036580863548334)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $rowValues[$j]
    }
}

Write-Host "Updated $($data.Count) rows of vm_pu results."
